$d = $word.ActiveDocument

# 1) Remove the "Meta description" paragraph that directly follows the
#    title heading at the top of the document.
$metaPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Meta description*") {
        $metaPara = $candidate
        break
    }
}
if ($metaPara -ne $null) {
    [void]$metaPara.Range.Delete()
}

# 2) Replace the closing "DALLE feature-image prompt" paragraph (now the
#    last paragraph in the document) with two new paragraphs: a bold
#    restatement of the page title, followed by the (now italic) meta
#    description sentence that used to live at the top of the document.
$lastPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "Create a feature image fit for D'Cirque!*") {
        $lastPara = $candidate
        break
    }
}
if ($lastPara -ne $null) {
    $target = $d.Range($lastPara.Range.Start, $lastPara.Range.End)

    $fragment = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
        '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play D''Cirque for Free - Review of Circus-Themed Slot Game</w:t></w:r></w:p>' + `
        '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of D''Cirque, a circus-themed slot game with superb graphics and a fair RTP of 96.11%. Play for free and win big with expanding Wilds and Free Spins.</w:t></w:r></w:p>' + `
        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    [void]$target.InsertXML($fragment)
}
